$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 1 through 6, keeping only the last data row
# (FC Barcelona vs Villareal CF / May 22nd 2022 / 84), which shifts up to row 1
$ws.Range("A1:C6").EntireRow.Delete()

# Update the score value for the remaining row from 84 to 62, keeping it
# stored as text (shared string) like the original cell, then drop the
# temporary text number-format so no stray style is left behind.
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "62"
$ws.Range("C1").ClearFormats()
